# "Generate Report for Archive"
#
# The handoff status text changes from "Ready for handoff" to
# "In Translation" everywhere it appears (the shared string is reused
# across the Overview sheet's E2/F2 "status" cells and the C2 "Status"
# cell on each per-language sheet). Because Excel auto-sizes those status
# columns to fit their content, the now-shorter text also narrows the
# columns that display it.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: E2 (zh-cn status) and F2 (de-de status) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: C2 status ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: C2 status ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Columns.Item(3).ColumnWidth = 12.5
